$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells in row 1 (AD1, AE1) ---
# Copy style (yellow fill header style) from existing header cell AC1
$null = $ws.Range("AC1").Copy()
$null = $ws.Range("AD1").PasteSpecial(-4122)
$null = $ws.Range("AC1").Copy()
$null = $ws.Range("AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths for new columns AD (30) and AE (31) ---
# (target character widths are 70.42578125 / 16.28515625; the interop
# runtime quantizes ColumnWidth to 1/6-character steps, so we use the
# input value that lands closest on the nearest achievable width)
$ws.Range("AD1").EntireColumn.ColumnWidth = 69.66666667
$ws.Range("AE1").EntireColumn.ColumnWidth = 15.5

# Values are entered in this precise order so that the shared-strings table
# is built up in the same sequence as the source workbook.
$ws.Range("A18").Value = "HeaderLinks"
$ws.Range("AD1").Value = "HeaderNames"
$ws.Range("AD19").Value = "Dryers,Straighteners,Curling Irons,Hair Brushes & Elastics,Specialty"
$ws.Range("A19").Value = "HeaderMobileLinks"
$ws.Range("AD18").Value = "Collections,Dryers,Straighteners,Curling Irons,Hair Brushes & Elastics,Specialty"
$ws.Range("AE1").Value = "Promocode"
$ws.Range("A20").Value = "Promocode"
$ws.Range("AE20").Value = "20OFF!"
$ws.Range("S20").Value = "1"" Long Lasting Curls Heated Silicone Brush"

# --- Sheet view changes: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 120
$null = $ws.Range("X5").Select()
